$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1 (06:00 -> 06:15)
$ws.Range("F1").Value = "Last status check on: 09.02.2022 06:15"

# D4: was inline text "+0.4", now a plain numeric value 0.4
$ws.Range("D4").Value = 0.4

# E4: was inline text "2022-02-09 06:00:17", now a numeric date serial
# formatted the same way as the other date cells in column E (style s="2").
$ws.Range("E4").Value = 44601.25019675926
$ws.Range("E4").NumberFormat = $ws.Range("E2").NumberFormat
